# Implement DuckDuckGo search engine test cases.
#
# Sheet "Google" (Worksheets(1)): the old keyword "no such element exception
# selenium" in A2 is replaced with "automation testing", and two new
# keywords are appended in A3 / A4.
#
# Sheet "DuckDuckGo" (Worksheets(2)): keeps its existing keywords, but
# becomes the active/selected sheet with a new selection.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Google
$ws2 = $wb.Worksheets.Item(2)   # DuckDuckGo

# Clear the cells whose text is being replaced/reused so the shared-string
# table doesn't keep stale, unreferenced entries around.
$ws1.Range("A2").ClearContents()
$ws2.Range("A3").ClearContents()
$ws2.Range("A4").ClearContents()

# Google sheet: new / updated keywords.
$ws1.Range("A2").Value = "automation testing"
$ws1.Range("A3").Value = "quality assurance"
$ws1.Range("A4").Value = "ruby is the best programming language"

# DuckDuckGo sheet: re-enter its existing keywords.
$ws2.Range("A3").Value = "privacy"
$ws2.Range("A4").Value = "hire me please"

# Leave the Google sheet's remembered selection on A3.
$null = $ws1.Activate()
$null = $ws1.Range("A3").Select()

# DuckDuckGo becomes the active tab, with B7 selected.
$null = $ws2.Activate()
$null = $ws2.Range("B7").Select()
